$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Item Name (column D) and UOM (column E) values, rows 2-27, reordered per
# the target workbook state (sharedStrings reshuffled; cell text follows).
$ws.Range("D2").Value = "Desodin 60ml Syrup"
$ws.Range("D3").Value = "Dinafex 120mg Tablet"
$ws.Range("D4").Value = "Dinafex 60mg Tablet"
$ws.Range("D5").Value = "Dinafex 180mg Tablet"
$ws.Range("D6").Value = "Dorenta 50mg Tablet"
$ws.Range("D7").Value = "Etorix 60mg Tablet - 40's"
$ws.Range("D8").Value = "Etorix 90mg Tablet"
$ws.Range("D9").Value = "Etorix 120mg Tablet"
$ws.Range("D10").Value = "Fenobac 100ml Syrup"
$ws.Range("D11").Value = "Flucloxin 500mg Capsule - 36's"
$ws.Range("D12").Value = "Flucloxin 500mg Capsule"
$ws.Range("D13").Value = "Geminox 320mg Tablet - 8's"
$ws.Range("D14").Value = "Ketonic 30mg IM/IV Injection - 4's"
$ws.Range("D15").Value = "Ketonic 30mg Injection"
$ws.Range("D16").Value = "Ketonic 10mg Tablet"
$ws.Range("D17").Value = "Kynol TR 100mg Capsule"
$ws.Range("D18").Value = "Kynol D 25mg Tablet"
$ws.Range("D19").Value = "Kynol TR 200mg Capsule"
$ws.Range("D20").Value = "Naprox Plus 500mg Tablet - 30's"
$ws.Range("D21").Value = "Oradin Plus Tablet - 40's"
$ws.Range("D22").Value = "Osticare Tablet 24's"
$ws.Range("D23").Value = "Sk-Mox 500mg Capsule"
$ws.Range("D24").Value = "Zithrox 250mg Tablet - 6's"
$ws.Range("D25").Value = "Zithrox 15ml Suspension"
$ws.Range("D26").Value = "Zithrox 30ml Dry Suspension"
$ws.Range("D27").Value = "Zithrox 500mg Tablet"

$ws.Range("E2").Value = "60 ml"
$ws.Range("E3").Value = "30's"
$ws.Range("E4").Value = "30's"
$ws.Range("E5").Value = "30's"
$ws.Range("E6").Value = "50's"
$ws.Range("E7").Value = "40's"
$ws.Range("E8").Value = "30's"
$ws.Range("E9").Value = "20's"
$ws.Range("E10").Value = "100ml"
$ws.Range("E11").Value = "36 's"
$ws.Range("E12").Value = "30 's"
$ws.Range("E13").Value = "8 's"
$ws.Range("E14").Value = "4's"
$ws.Range("E15").Value = "5 's"
$ws.Range("E16").Value = "20's"
$ws.Range("E17").Value = "50 's"
$ws.Range("E18").Value = "60 's"
$ws.Range("E19").Value = "30 's"
$ws.Range("E20").Value = "30 's"
$ws.Range("E21").Value = "40 's"
$ws.Range("E22").Value = "24's"
$ws.Range("E23").Value = "48 's"
$ws.Range("E24").Value = "6's"
$ws.Range("E25").Value = "15 ml"
$ws.Range("E26").Value = "30ml"
$ws.Range("E27").Value = "6 's"
